$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.773.54"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "3.670.36"
$ws.Range("D4").Value = "'2.71"
$ws.Range("E4").Value = "  +42.55%  "
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "'227.83"
$ws.Range("E6").Value = "  -4.69%  "
$ws.Range("D7").Value = "'646.63"
$ws.Range("E7").Value = "  -2.42%  "
$ws.Range("D8").Value = "'0.425"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "'1.19"
$ws.Range("E9").Value = "  +11.33%  "
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").Value = "3.667.00"
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").Value = "'49.69"
$ws.Range("E12").Value = "  +11.03%  "
$ws.Range("D13").Value = "'0.212"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D14").Value = "'0.0000294"
$ws.Range("E14").Value = "  -8.51%  "
$ws.Range("D15").Value = "'6.69"
$ws.Range("E15").Value = "  -3.29%  "
$ws.Range("D16").Value = "4.385.76"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").Value = "96.663.30"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "'21.58"
$ws.Range("E18").Value = "  +13.78%  "
$ws.Range("D19").Value = "'8.92"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").Value = "'14.19"
$ws.Range("E20").Value = "  +7.53%  "
$ws.Range("D21").Value = "3.657.80"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "'0.543"
$ws.Range("E22").Value = "  +7.81%  "
$ws.Range("D23").Value = "'0.277"
$ws.Range("E23").Value = "  +43.70%  "
$ws.Range("D24").Value = "'532.87"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").Value = "'3.28"
$ws.Range("E25").Value = "  -5.88%  "
$ws.Range("D26").Value = "'124.73"
$ws.Range("E26").Value = "  +14.78%  "
$ws.Range("D27").Value = "'0.0000203"
$ws.Range("E27").Value = "  -10.24%  "
$ws.Range("D28").Value = "'6.86"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.848.24"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "'13.04"
$ws.Range("E30").Value = "  -5.24%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'13.21"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'3.02"
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("B34").Value = "Cronos"
$ws.Range("C34").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D34").Value = "'0.184"
$ws.Range("E34").Value = "  -5.25%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'33.29"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").Value = "'0.625"
$ws.Range("E36").Value = "  +4.35%  "
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'1.78"
$ws.Range("E38").Value = "  -4.01%  "
$ws.Range("D39").Value = "'605.99"
$ws.Range("E39").Value = "  -6.77%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "'8.54"
$ws.Range("E41").Value = "  -3.75%  "
$ws.Range("D42").Value = "'7.15"
$ws.Range("E42").Value = "  +3.55%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'42.41"
$ws.Range("E43").Value = "  +1.90%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.502"
$ws.Range("E44").Value = "  +4.65%  "
$ws.Range("D45").Value = "'0.0503"
$ws.Range("E45").Value = "  +8.96%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").Value = "'0.160"
$ws.Range("E46").Value = "  -4.82%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'0.971"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("B48").Value = "ImmutableX"
$ws.Range("C48").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D48").Value = "'1.96"
$ws.Range("E48").Value = "  -4.36%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'234.01"
$ws.Range("E49").Value = "  +12.63%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.32"
$ws.Range("E50").Value = "  -4.66%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'8.85"
$ws.Range("E51").Value = "  +0.31%  "
